# McCarthy discussion and minor correction up to GAMLSS section
# Row 10 in the "Hourly" worksheet is the McCarthy reference row. This script
# updates the " Length" (D), the numeric reference count (E), the "Metric"
# (G) and the "Probabilistic" (H) columns for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hourly")

$ws.Range("G10").Value = "Calibration"
$ws.Range("D10").Value = "n/a"
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = "Yes"

# Restore the cursor/selection to where the author left it when saving.
$ws.Range("H10").Select()
